# Update the "Syllabi Scores" table: rows were re-ordered / re-mapped to new
# course codes (some codes were also cleaned up / simplified), while the
# score data (columns B:W) for each course travels with its original row.
# Row 1 (header) and row 12 (BAEN 540) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Syllabi Scores")

# Snapshot of the current data block (rows 2-42, columns A-W) before any
# writes happen, so that source rows are not clobbered while building the
# re-ordered table.
$origRange = $ws.Range("A2:W42")
$orig = $origRange.Value2

# Mapping: for every row in the NEW layout, which OLD row (2-42) supplies
# the B:W score values, and what the NEW course-code label (column A)
# should be. An empty Label means "keep the existing label" (used for the
# single untouched row, 12).
$rowMap = @(
    @{ New = 2; Old = 16; Label = "BAEN 580" },
    @{ New = 3; Old = 19; Label = "BAPA 580" },
    @{ New = 4; Old = 18; Label = "BAEN 580C" },
    @{ New = 5; Old = 17; Label = "BAEN 580A" },
    @{ New = 6; Old = 7; Label = "BA 562" },
    @{ New = 7; Old = 6; Label = "BA 560" },
    @{ New = 8; Old = 15; Label = "BAEN550B" },
    @{ New = 9; Old = 14; Label = "BAEN550A" },
    @{ New = 10; Old = 13; Label = "BAEN 549" },
    @{ New = 11; Old = 5; Label = "APSC 540" },
    @{ New = 12; Old = 12; Label = "" },
    @{ New = 13; Old = 20; Label = "BASM 516" },
    @{ New = 14; Old = 11; Label = "BAEN 509" },
    @{ New = 15; Old = 10; Label = "BAEN 506" },
    @{ New = 16; Old = 9; Label = "BAEN 505" },
    @{ New = 17; Old = 8; Label = "BAEN 502" },
    @{ New = 18; Old = 41; Label = "FCOR 502" },
    @{ New = 19; Old = 36; Label = "COMM 497" },
    @{ New = 20; Old = 35; Label = "COMM 489" },
    @{ New = 21; Old = 4; Label = "APSC 486" },
    @{ New = 22; Old = 34; Label = "COMM 486G" },
    @{ New = 23; Old = 33; Label = "COMM 485" },
    @{ New = 24; Old = 32; Label = "COMM 482" },
    @{ New = 25; Old = 31; Label = "COMM 470" },
    @{ New = 26; Old = 30; Label = "COMM 466" },
    @{ New = 27; Old = 40; Label = "ENPH 459" },
    @{ New = 28; Old = 39; Label = "DES 445" },
    @{ New = 29; Old = 3; Label = "APSC 440" },
    @{ New = 30; Old = 21; Label = "BMEG 401" },
    @{ New = 31; Old = 29; Label = "COMM 389" },
    @{ New = 32; Old = 37; Label = "COMM 388" },
    @{ New = 33; Old = 28; Label = "COMM 387" },
    @{ New = 34; Old = 27; Label = "COMM 386L" },
    @{ New = 35; Old = 26; Label = "COMM 386I" },
    @{ New = 36; Old = 25; Label = "COMM 383" },
    @{ New = 37; Old = 2; Label = "APSC 383" },
    @{ New = 38; Old = 24; Label = "COMM 382" },
    @{ New = 39; Old = 42; Label = "IGEN 340" },
    @{ New = 40; Old = 38; Label = "COMR 280" },
    @{ New = 41; Old = 22; Label = "COMM 280B" },
    @{ New = 42; Old = 23; Label = "COMM 280A" }
)

$numRows = 41
$numCols = 23
$newArr = New-Object 'object[,]' $numRows,$numCols

foreach ($entry in $rowMap) {
    $newRowIdx = $entry.New - 2   # 0-based row offset into $newArr
    $oldRowIdx = $entry.Old - 2   # 0-based row offset into $orig (also 0-based after Value2? use 1-based accessor below)

    for ($col = 0; $col -lt $numCols; $col++) {
        $newArr[$newRowIdx, $col] = $orig[$entry.Old - 1, $col + 1]
    }

    if ($entry.Label -ne "") {
        $newArr[$newRowIdx, 0] = $entry.Label
    }
}

$ws.Range("A2:W42").Value2 = $newArr
